$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 24
$ws.Range("D2").Value = 20
$ws.Range("C3").Value = 41
$ws.Range("D3").Value = 32
$ws.Range("C4").Value = 82
$ws.Range("D4").Value = 62
$ws.Range("C5").Value = 70
$ws.Range("D5").Value = 63
$ws.Range("D6").Value = 80
$ws.Range("C7").Value = 70
$ws.Range("D7").Value = 53
$ws.Range("C8").Value = 45
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = 41
$ws.Range("D10").Value = 21
$ws.Range("C11").Value = 61
$ws.Range("D11").Value = 51
$ws.Range("D12").Value = 30
$ws.Range("C14").Value = 91
$ws.Range("D14").Value = 77
$ws.Range("C15").Value = 55
$ws.Range("D15").Value = 43
$ws.Range("C17").Value = 32
$ws.Range("D17").Value = 27
$ws.Range("C19").Value = 45
$ws.Range("D19").Value = 36
$ws.Range("C20").Value = 30
$ws.Range("D20").Value = 23
$ws.Range("D21").Value = 45
$ws.Range("D23").Value = 16
$ws.Range("D24").Value = 54
$ws.Range("C26").Value = 37
$ws.Range("D26").Value = 29
$ws.Range("C27").Value = 55
$ws.Range("D27").Value = 41
$ws.Range("C28").Value = 54
$ws.Range("D28").Value = 40
$ws.Range("D29").Value = 46
$ws.Range("C30").Value = 51
$ws.Range("D30").Value = 36
$ws.Range("C31").Value = 57
$ws.Range("D31").Value = 45
$ws.Range("C32").Value = 43
$ws.Range("D32").Value = 32
$ws.Range("C33").Value = 48
$ws.Range("C35").Value = 53
$ws.Range("D35").Value = 45
$ws.Range("C37").Value = 55
$ws.Range("D37").Value = 35
$ws.Range("C38").Value = 63
$ws.Range("D38").Value = 48
$ws.Range("C39").Value = 61
$ws.Range("D39").Value = 53
$ws.Range("D41").Value = 43
$ws.Range("C42").Value = 73
$ws.Range("D42").Value = 55
$ws.Range("C43").Value = 61
$ws.Range("D43").Value = 48
$ws.Range("C44").Value = 64
$ws.Range("D44").Value = 53
$ws.Range("C45").Value = 36
$ws.Range("D45").Value = 22
$ws.Range("C46").Value = 65
$ws.Range("D46").Value = 47
$ws.Range("C47").Value = 48
$ws.Range("D47").Value = 41
$ws.Range("D48").Value = 28
$ws.Range("C49").Value = 57
$ws.Range("D49").Value = 46
$ws.Range("D51").Value = 29
$ws.Range("D52").Value = 36
$ws.Range("C53").Value = 82
$ws.Range("D53").Value = 71
$ws.Range("C54").Value = 51
$ws.Range("D54").Value = 30
$ws.Range("C55").Value = 64
$ws.Range("D55").Value = 44
$ws.Range("D57").Value = 40
$ws.Range("C58").Value = 66
$ws.Range("D58").Value = 55
$ws.Range("C60").Value = 53
$ws.Range("C61").Value = 32
$ws.Range("D61").Value = 24
$ws.Range("C62").Value = 55
$ws.Range("D62").Value = 40
$ws.Range("D63").Value = 93
$ws.Range("C65").Value = 30
$ws.Range("D65").Value = 26
$ws.Range("D68").Value = 34
$ws.Range("C69").Value = 95
$ws.Range("C72").Value = 39
$ws.Range("D72").Value = 20
$ws.Range("D73").Value = 36
$ws.Range("D74").Value = 93
$ws.Range("C76").Value = 28
$ws.Range("C78").Value = 75
$ws.Range("D78").Value = 57
$ws.Range("C79").Value = 108
$ws.Range("D79").Value = 108
$ws.Range("D80").Value = 89
$ws.Range("D81").Value = 49
$ws.Range("C82").Value = 26
$ws.Range("D82").Value = 15
$ws.Range("D83").Value = 60
$ws.Range("C84").Value = 185
$ws.Range("D84").Value = 119
$ws.Range("C92").Value = 223
$ws.Range("D92").Value = 146
$ws.Range("C93").Value = 5179
$ws.Range("D93").Value = 3996
